$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as captured by the scraper run.
# Column D values that render as plain decimal numbers must be forced to Text format
# so Excel keeps them as the scraped string (matching the sheet convention of storing
# "Price" as text) instead of silently converting them to floating point numbers.

$ws.Range('D2').Value = '42.795.50'
$ws.Range('E2').Value = '  -0.64%  '
$ws.Range('D3').Value = '2.293.36'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.89'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.64'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.508'
$ws.Range('E7').Value = '  -2.04%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -2.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.61'
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.37'
$ws.Range('E12').Value = '  +3.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.119'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.74'
$ws.Range('E14').Value = '  -1.91%  '
$ws.Range('D15').Value = '2.651.78'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').Value = '2.294.56'
$ws.Range('E16').Value = '  -0.39%  '
$ws.Range('E17').Value = '  -0.96%  '
$ws.Range('D18').Value = '42.722.48'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.03'
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').Value = '0.0₃0899'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.20'
$ws.Range('E22').Value = '  -1.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.17'
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.17'
$ws.Range('E28').Value = '  +0.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.58'
$ws.Range('E29').Value = '  +2.03%  '
$ws.Range('E30').Value = '  +0.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.07'
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.28'
$ws.Range('E32').Value = '  +0.90%  '
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.76'
$ws.Range('E34').Value = '  +0.22%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.98'
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.74'
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0691'
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('E39').Value = '  -0.70%  '
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('E41').Value = '  -1.41%  '
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('D43').Value = '2.000.34'
$ws.Range('E43').Value = '  -0.45%  '
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.14'
$ws.Range('E45').Value = '  +4.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.97'
$ws.Range('E46').Value = '  -2.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.08'
$ws.Range('E47').Value = '  -5.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  -2.72%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.88'
$ws.Range('E49').Value = '  +5.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '53.76'
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('D51').Value = '2.520.86'
$ws.Range('E51').Value = '  -0.36%  '
